$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2304760901931839
$ws.Range("C2").Value = 0.6389819110794862
$ws.Range("D2").Value = 0.5920202495404218
$ws.Range("E2").Value = 0.769428521397811
$ws.Range("F2").Value = 0.7618103084836596
$ws.Range("G2").Value = 14
